$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: add the new Mongo/Flask note in I9, with wrapping, taller row ---
$ws.Rows.Item(9).RowHeight = 34
$ws.Range("I9").Value = "Mongo는 프로시저가 없어`nFlask에서 설정"
$ws.Range("I9").WrapText = $true

# --- Rows 10 & 11: swap the DevOps / FE task entries ---
# Row 10 becomes the (now complete) DevOps / Docker Flask 세팅 task
$ws.Range("E10").Value = $true
$ws.Range("F10").Value = "DevOps"
$ws.Range("G10").Value = "Docker Flask 세팅"
$ws.Range("K10").Value = "2026 09 26"
$ws.Range("L10").Value = "2026 09 26"

# Row 11 becomes the (still open) FE / 데이터 바인딩 task
$ws.Range("F11").Value = "FE"
$ws.Range("G11").Value = "데이터 바인딩"

# L11 loses its outer (medium) right border, matching L5's plain style
$ws.Range("L5").Copy()
$ws.Range("L11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Selection moves to J10 ---
[void]$ws.Range("J10").Select()
